# Edit script applying the commit "debugging system, now sysA matches results
# from Trimmer et al." to the LCA report workbook.
#
# Summary of the change:
#  - Header label "Total GlobalWarming Ratio" -> "Category GlobalWarming Ratio"
#    (appears at F1, F30, D36, D48)
#  - Header label "WasteStream" -> "Stream" (A36)
#  - SanUnit numbering shifted by one for the Concrete/Excavation/Plastic/Steel
#    breakdowns: A4->A5, A5->A6, A6->A7, A7->A8
#  - All of the LCA quantity/ratio figures were recomputed with the corrected
#    system (sysA) numbers
#  - The now-unused helper column H (=E/40000/8 utilisation calc) is removed
#  - Column D and F widths widened from 22 to 25 to fit the new header text
#  - Selection moved from H2:H27 to E43:E46

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LCA")

# --- Header / label text updates -----------------------------------------
$ws.Range("F1").Value = "Category GlobalWarming Ratio"
$ws.Range("F30").Value = "Category GlobalWarming Ratio"
$ws.Range("D36").Value = "Category GlobalWarming Ratio"
$ws.Range("D48").Value = "Category GlobalWarming Ratio"
$ws.Range("A36").Value = "Stream"

# --- SanUnit code relabeling (A4..A7 shift to A5..A8) ---------------------
$ws.Range("B6").Value = "A5"
$ws.Range("B10").Value = "A5"
$ws.Range("B7").Value = "A8"
$ws.Range("B23").Value = "A8"
$ws.Range("B11").Value = "A6"
$ws.Range("B17").Value = "A6"
$ws.Range("B12").Value = "A7"
$ws.Range("B18").Value = "A7"

# --- Recomputed LCA figures ------------------------------------------------
$ws.Range("C2").Value = 5906254.725
$ws.Range("E2").Value = 1653751.323
$ws.Range("F2").Value = 0.05746078661061285
$ws.Range("C3").Value = 5906254.725
$ws.Range("E3").Value = 1653751.323
$ws.Range("F3").Value = 0.05746078661061285
$ws.Range("C4").Value = 18229181.25
$ws.Range("E4").Value = 19687515.75
$ws.Range("F4").Value = 0.6840569834596768
$ws.Range("C5").Value = 18229181.25
$ws.Range("E5").Value = 19687515.75
$ws.Range("F5").Value = 0.6840569834596768
$ws.Range("F6").Value = 0.004513933361288044
$ws.Range("F7").Value = 0.06314833765316082
$ws.Range("F8").Value = 0.06766227101444887
$ws.Range("C9").Value = 95208.40950000001
$ws.Range("D9").Value = 0.7124265141203362
$ws.Range("E9").Value = 50460.45703500001
$ws.Range("F9").Value = 0.001753285100145697
$ws.Range("D10").Value = 0.01085914278195631
$ws.Range("F10").Value = 0.00002672440295609602
$ws.Range("D11").Value = 0.1041607262282339
$ws.Range("F11").Value = 0.0002563400514954312
$ws.Range("D12").Value = 0.1725536168694736
$ws.Range("F12").Value = 0.00042465528645723
$ws.Range("C13").Value = 133639.6212282441
$ws.Range("E13").Value = 70828.99925096935
$ws.Range("F13").Value = 0.002461004841054454
$ws.Range("C14").Value = 33333360
$ws.Range("E14").Value = 500000.4
$ws.Range("F14").Value = 0.01737287577040449
$ws.Range("C15").Value = 33333360
$ws.Range("E15").Value = 500000.4
$ws.Range("F15").Value = 0.01737287577040449
$ws.Range("C16").Value = 262500.21
$ws.Range("D16").Value = 0.9425426423914033
$ws.Range("E16").Value = 517125.4137
$ws.Range("F16").Value = 0.01796789676554084
$ws.Range("D17").Value = 0.01630139246494967
$ws.Range("F17").Value = 0.0003107570138170469
$ws.Range("D18").Value = 0.04115596514364706
$ws.Range("F18").Value = 0.0007845651748031642
$ws.Range("C19").Value = 278502.2111402715
$ws.Range("E19").Value = 548649.3559463349
$ws.Range("F19").Value = 0.01906321895416106
$ws.Range("C20").Value = 82614649.425
$ws.Range("E20").Value = 991375.7931
$ws.Range("F20").Value = 0.0344460694337695
$ws.Range("C21").Value = 82614649.425
$ws.Range("E21").Value = 991375.7931
$ws.Range("F21").Value = 0.0344460694337695
$ws.Range("C22").Value = 874349.6578125
$ws.Range("D22").Value = 0.9265632219489316
$ws.Range("E22").Value = 2229591.627421875
$ws.Range("F22").Value = 0.07746877474884861
$ws.Range("D23").Value = 0.07343677805106841
$ws.Range("F23").Value = 0.00613995578753175
$ws.Range("C24").Value = 943648.1365765785
$ws.Range("E24").Value = 2406302.748270275
$ws.Range("F24").Value = 0.08360873053638036
$ws.Range("C25").Value = 4947.920625
$ws.Range("E25").Value = 974740.3631249999
$ws.Range("F25").Value = 0.03386805937949167
$ws.Range("C26").Value = 4947.920625
$ws.Range("E26").Value = 974740.3631249999
$ws.Range("F26").Value = 0.03386805937949167
$ws.Range("E27").Value = 28780520.08244796
$ws.Range("C31").Value = 4500003.600000004
$ws.Range("E31").Value = 873000.6984000008
$ws.Range("C32").Value = 4500003.600000004
$ws.Range("E32").Value = 873000.6984000008
$ws.Range("E33").Value = 873000.6984000008
$ws.Range("B37").Value = 2940576.254688381
$ws.Range("C37").Value = 82336135.13127466
$ws.Range("D37").Value = 1.374056830128091
$ws.Range("B38").Value = 653.0315412936889
$ws.Range("C38").Value = 173053.3584428275
$ws.Range("D38").Value = 0.002887980456768505
$ws.Range("B39").Value = 1806400.205981089
$ws.Range("C39").Value = -2709600.308971634
$ws.Range("D39").Value = -0.04521884353113698
$ws.Range("B40").Value = 3024555.141536747
$ws.Range("C40").Value = -16332597.76429843
$ws.Range("D40").Value = -0.2725646215478565
$ws.Range("B41").Value = 315521.7839492012
$ws.Range("C41").Value = -1546056.741351086
$ws.Range("D41").Value = -0.02580118464185861
$ws.Range("B42").Value = 275430.4794482073
$ws.Range("C42").Value = -413145.719172311
$ws.Range("D42").Value = -0.00689473335567418
$ws.Range("B43").Value = 116754.6613344789
$ws.Range("C43").Value = -630475.171206186
$ws.Range("D43").Value = -0.01052161015137299
$ws.Range("B44").Value = 194976.3844279827
$ws.Range("C44").Value = -955384.2836971155
$ws.Range("D44").Value = -0.01594381735696042
$ws.Range("C45").Value = 59921928.50102073
$ws.Range("B49").Value = 456960
$ws.Range("C49").Value = 68544
$ws.Range("C50").Value = 68544

# --- Remove the now-unused helper column H (=E/40000/8) --------------------
$ws.Range("H1:H27").ClearContents()

# --- Column width adjustments for the longer header text -------------------
# (ColumnWidth is in "characters"; the stored XML width includes ~5/6 chars of
# padding, so subtract that to land on an XML-stored width of exactly 25)
$ws.Range("D1").EntireColumn.ColumnWidth = 24.166666666666668
$ws.Range("F1").EntireColumn.ColumnWidth = 24.166666666666668

# --- Update the active selection -------------------------------------------
$ws.Range("E43:E46").Select()
